$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Shrimp")

# Row 22: SeaFood Bar&Market -> Mango Mama
$ws.Range("A22").Value = "Mango Mama"
$ws.Range("B22").Value = "Świętego Mikołaja 18, 50-128 Wrocław "
$ws.Range("C22").Value = "535 653 692 "
$ws.Range("D22").Value = "mangomama.pl "

# Row 23: Shrimp Nature -> Wok in
$ws.Range("A23").Value = "Wok in"
$ws.Range("B23").Value = "Sukiennice 1/2, 50-116 Wrocław "
$ws.Range("C23").Value = "781 026 510 "
$ws.Range("D23").Value = "wokin.pl "
$ws.Range("E23").Value = " 4,6"

# Row 24: Sea-Food Crab Meat Sokolnicza 7/17 Wrocław -> Vertigo Jazz Club & Restaurant
$ws.Range("A24").Value = "Vertigo Jazz Club & Restaurant"
$ws.Range("B24").Value = "Oławska 13, 50-123 Wrocław "
$ws.Range("C24").Value = "71 718 25 81 "
$ws.Range("D24").Value = "vertigojazz.pl "
$ws.Range("E24").Value = " 4,7"
